$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.718.89"
$ws.Range("E2").Value = "  +0.24%  "

# Row 3
$ws.Range("E3").Value = "  +0.10%  "

# Row 4
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.29%  "

# Row 5
$ws.Range("D5").Value = "211.39"
$ws.Range("E5").Value = "  -0.23%  "

# Row 6
$ws.Range("D6").Value = "0.512"
$ws.Range("E6").Value = "  -0.24%  "

# Row 7
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.22%  "

# Row 8
$ws.Range("D8").Value = "0.0618"
$ws.Range("E8").Value = "  -0.04%  "

# Row 9
$ws.Range("E9").Value = "  +0.53%  "

# Row 10
$ws.Range("D10").Value = "19.63"
$ws.Range("E10").Value = "  +0.41%  "

# Row 11
$ws.Range("D11").Value = "0.0842"
$ws.Range("E11").Value = "  +0.73%  "

# Row 12
$ws.Range("D12").Value = "1.826.51"
$ws.Range("E12").Value = "  +0.12%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "4.05"
$ws.Range("E13").Value = "  +0.22%  "

# Row 14
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.565.01"
$ws.Range("E14").Value = "  -1.67%  "

# Row 15
$ws.Range("D15").Value = "0.523"
$ws.Range("E15").Value = "  -0.19%  "

# Row 16
$ws.Range("D16").Value = "65.36"
$ws.Range("E16").Value = "  +0.98%  "

# Row 17
$ws.Range("D17").Value = "26.688.87"
$ws.Range("E17").Value = "  +0.12%  "

# Row 18
$ws.Range("D18").Value = "0.0₃0748"
$ws.Range("E18").Value = "  +1.98%  "

# Row 19
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "209.78"
$ws.Range("E19").Value = "  +0.06%  "

# Row 20
$ws.Range("E20").Value = "  +0.27%  "

# Row 21
$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").Value = "7.18"
$ws.Range("E21").Value = "  +0.82%  "

# Row 22
$ws.Range("E22").Value = "  +0.46%  "

# Row 23
$ws.Range("D23").Value = "2.30"
$ws.Range("E23").Value = "  -0.29%  "

# Row 24
$ws.Range("D24").Value = "8.96"
$ws.Range("E24").Value = "  +0.16%  "

# Row 25
$ws.Range("D25").Value = "142.57"
$ws.Range("E25").Value = "  -1.86%  "

# Row 26
$ws.Range("E26").Value = "  +0.31%  "

# Row 27
$ws.Range("D27").Value = "7.11"
$ws.Range("E27").Value = "  -0.58%  "

# Row 28
$ws.Range("D28").Value = "0.114"
$ws.Range("E28").Value = "  -0.99%  "

# Row 29
$ws.Range("D29").Value = "15.43"
$ws.Range("E29").Value = "  +0.91%  "

# Row 30
$ws.Range("D30").Value = "0.0518"
$ws.Range("E30").Value = "  +1.93%  "

# Row 31
$ws.Range("E31").Value = "  -0.39%  "

# Row 32
$ws.Range("E32").Value = "  +0.62%  "

# Row 33
$ws.Range("D33").Value = "2.97"
$ws.Range("E33").Value = "  +1.10%  "

# Row 34
$ws.Range("D34").Value = "1.305.12"
$ws.Range("E34").Value = "  +1.54%  "

# Row 35
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "2.47"
$ws.Range("E35").Value = "  +1.11%  "

# Row 36
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "0.614"
$ws.Range("E36").Value = "  -3.83%  "

# Row 37
$ws.Range("E37").Value = "  +0.14%  "

# Row 38
$ws.Range("E38").Value = "  -0.31%  "

# Row 39
$ws.Range("D39").Value = "1.11"
$ws.Range("E39").Value = "  +20.41%  "

# Row 40
$ws.Range("D40").Value = "0.823"
$ws.Range("E40").Value = "  -2.96%  "

# Row 41
$ws.Range("D41").Value = "5.43"
$ws.Range("E41").Value = "  -1.18%  "

# Row 42
$ws.Range("E42").Value = "  +0.22%  "

# Row 43
$ws.Range("D43").Value = "0.782"
$ws.Range("E43").Value = "  -0.69%  "

# Row 44
$ws.Range("D44").Value = "63.21"
$ws.Range("E44").Value = "  -2.10%  "

# Row 45
$ws.Range("D45").Value = "1.737.94"
$ws.Range("E45").Value = "  +0.09%  "

# Row 46
$ws.Range("D46").Value = "91.26"
$ws.Range("E46").Value = "  +1.16%  "

# Row 47
$ws.Range("D47").Value = "1.57"
$ws.Range("E47").Value = "  -2.15%  "

# Row 48
$ws.Range("D48").Value = "0.0₆0106"
$ws.Range("E48").Value = "  -0.41%  "

# Row 49
$ws.Range("E49").Value = "  -1.83%  "

# Row 50
$ws.Range("D50").Value = "0.0516"
$ws.Range("E50").Value = "  +1.40%  "

# Row 51
$ws.Range("E51").Value = "  +0.25%  "
